$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'26.517.18"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "'1.808.33"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'308.36"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("D7").Value = "'0.4560"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("D8").Value = "'0.3665"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "'0.07137"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "'0.8802"
$ws.Range("D11").Value = "'0.07755"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'19.40"
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "'1.806.95"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'5.281"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "'6.372"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "'86.17"
$ws.Range("E16").Value = "  -5.43%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'0.000008583"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'26.557.75"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").Value = "'14.26"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").Value = "'4.989"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "'10.44"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "'150.80"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'17.96"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'2.046"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'112.68"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").Value = "'4.848"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").Value = "'0.08676"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'3.043"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").Value = "'0.7315"
$ws.Range("E32").Value = "  -4.54%  "
$ws.Range("D33").Value = "'4.452"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").Value = "'1.116"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "'2.543"
$ws.Range("E36").Value = "  -7.70%  "
$ws.Range("D37").Value = "'1.081"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'0.05113"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'2.892"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'6.964"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "'0.5031"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "'0.1567"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Value = "'8.153"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "'1.008"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").Value = "'0.4615"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "'9.992"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "'100.75"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "'1.594"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "'0.05995"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "'64.06"
$ws.Range("E51").Value = "  -1.79%  "
